# Update "想去人数" (number of people interested) values in each sheet,
# reflecting the newer data snapshot generated at commit 2f43792.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1297
$ws1.Range("F3").Value = 2116
$ws1.Range("F4").Value = 433
$ws1.Range("F5").Value = 185
$ws1.Range("F6").Value = 436
$ws1.Range("F7").Value = 58
$ws1.Range("F8").Value = 540
$ws1.Range("F9").Value = 149
$ws1.Range("F10").Value = 94
$ws1.Range("F11").Value = 184
$ws1.Range("F12").Value = 823
$ws1.Range("F13").Value = 66
$ws1.Range("F15").Value = 4513
$ws1.Range("F16").Value = 2729
$ws1.Range("F17").Value = 858
$ws1.Range("F18").Value = 639
$ws1.Range("F19").Value = 339
$ws1.Range("F20").Value = 736
$ws1.Range("F21").Value = 1507
$ws1.Range("F22").Value = 56
$ws1.Range("F23").Value = 685
$ws1.Range("F24").Value = 298
$ws1.Range("F25").Value = 95
$ws1.Range("F26").Value = 216

# --- Sheet "演出" (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 7

# --- Sheet "本地生活" (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 129

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 129
$ws4.Range("F4").Value = 7
$ws4.Range("F5").Value = 1297
$ws4.Range("F6").Value = 2116
$ws4.Range("F7").Value = 433
$ws4.Range("F8").Value = 185
$ws4.Range("F9").Value = 436
$ws4.Range("F10").Value = 58
$ws4.Range("F11").Value = 540
$ws4.Range("F12").Value = 149
$ws4.Range("F13").Value = 94
$ws4.Range("F14").Value = 184
$ws4.Range("F15").Value = 823
$ws4.Range("F16").Value = 66
$ws4.Range("F20").Value = 4513
$ws4.Range("F21").Value = 2729
$ws4.Range("F22").Value = 858
$ws4.Range("F23").Value = 639
$ws4.Range("F24").Value = 339
$ws4.Range("F25").Value = 736
$ws4.Range("F26").Value = 1507
$ws4.Range("F27").Value = 56
$ws4.Range("F28").Value = 685
$ws4.Range("F29").Value = 298
$ws4.Range("F30").Value = 95
$ws4.Range("F31").Value = 216

$wb.Save()
